$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data rows (2-5) and the old extra header columns (I:K)
$ws.Rows("2:5").ClearContents()
$ws.Range("I1:K1").ClearContents()

# Rewrite the header row with the new, streamlined set of columns
$ws.Range("A1").Value = "Stock ID"
$ws.Range("B1").Value = "Price"
$ws.Range("C1").Value = "Yield"
$ws.Range("D1").Value = "Annual Yield"
$ws.Range("E1").Value = "`$price/annual yield"
$ws.Range("F1").Value = "Annual Yield for `$1k"
$ws.Range("G1").Value = "Updated:"
$ws.Range("H1").Value = "2019-07-30 08:50:26.785769"

# Column widths matching the streamlined layout
$ws.Columns("A").ColumnWidth = 7.1015625
$ws.Columns("B").ColumnWidth = 11.68359375
$ws.Columns("C").ColumnWidth = 10.68359375
$ws.Columns("D").ColumnWidth = 10.3125
$ws.Columns("E").ColumnWidth = 15.7890625
$ws.Columns("F").ColumnWidth = 16.41796875
$ws.Columns("G").ColumnWidth = 9.89453125
$ws.Columns("H").ColumnWidth = 23.9453125

# Move selection to row 2 (the first row to fill with data)
$ws.Range("A2:XFD2").Select()
